$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - this shifts the existing rows 6..16 down to 7..17,
# preserving all their data (matches the diff's "rotation" pattern).
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44949
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 2800
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = 2900
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Diguillín"
$ws.Range("S6").Value = 1450
$ws.Range("T6").Value = 2
